$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.00"
# or "605.92" are not auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.187.94"
$ws.Range("E2").Value = "  +3.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.621.80"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.92"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.91"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.620.80"
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +14.34%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.03"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("E15").Value = "  +8.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.55"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.170.28"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.619.12"
$ws.Range("E18").Value = "  +4.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "383.52"
$ws.Range("E19").Value = "  +8.86%  "
$ws.Range("E20").Value = "  +5.37%  "
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.25"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +8.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.757.72"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0961"
$ws.Range("E30").Value = "  +6.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "544.24"
$ws.Range("E31").Value = "  +5.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +3.06%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.74"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.17"
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  +7.22%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +8.72%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.61"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("E49").Value = "  +5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.532"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("E51").Value = "  +2.61%  "
